$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 71: group2 signature header
$ws.Range("B71").Value = "SmartRules String group2(String input1, String input2, String input3, Integer input4)"

# Row 72: column headers
$ws.Range("B72").Value = "input1"
$ws.Range("C72").Value = "input2"
$ws.Range("D72").Value = "input3"
$ws.Range("E72").Value = "input4"
$ws.Range("F72").Value = "Return"

# Row 73
$ws.Range("B73").Value = "A"
$ws.Range("D73").Value = "A"
$ws.Range("E73").Value = "1.. 3"
$ws.Range("F73").Value = "R1"

# Row 74
$ws.Range("B74").Value = "B"
$ws.Range("C74").Value = "A"
$ws.Range("D74").Value = "A"
$ws.Range("E74").Value = "4.. 6"
$ws.Range("F74").Value = "R2"

# Row 75
$ws.Range("B75").Value = "C"
$ws.Range("C75").Value = "C"
$ws.Range("D75").Value = "A"
$ws.Range("E75").Value = "7 ..9"
$ws.Range("F75").Value = "R3"

# Row 76
$ws.Range("C76").Value = "B"
$ws.Range("D76").Value = "B"
$ws.Range("E76").Value = "9 .. 10"
$ws.Range("F76").Value = "R4"

# Row 77
$ws.Range("D77").Value = "B"
$ws.Range("E77").Value = "< 22"
$ws.Range("F77").Value = "R5"

# Row 78
$ws.Range("E78").Value = "> 22"
$ws.Range("F78").Value = "R6"

# Row 80: test table title
$ws.Range("B80").Value = "Test group2"

# Row 81: description row
$ws.Range("B81").Value = "input1"
$ws.Range("C81").Value = "input2"
$ws.Range("D81").Value = "input3"
$ws.Range("E81").Value = "input4"
$ws.Range("F81").Value = "_res_"

# Row 82: column headers
$ws.Range("B82").Value = "input1"
$ws.Range("C82").Value = "input2"
$ws.Range("D82").Value = "input3"
$ws.Range("E82").Value = "input4"
$ws.Range("F82").Value = "_res_"

# Row 83
$ws.Range("B83").Value = "A"
$ws.Range("D83").Value = "A"
$ws.Range("E83").Value = 2
$ws.Range("F83").Value = "R1"

# Row 84
$ws.Range("B84").Value = "B"
$ws.Range("C84").Value = "A"
$ws.Range("D84").Value = "A"
$ws.Range("E84").Value = 5
$ws.Range("F84").Value = "R2"

# Row 85
$ws.Range("B85").Value = "C"
$ws.Range("C85").Value = "C"
$ws.Range("D85").Value = "A"
$ws.Range("E85").Value = 8
$ws.Range("F85").Value = "R3"

# Row 86
$ws.Range("C86").Value = "B"
$ws.Range("D86").Value = "B"
$ws.Range("E86").Value = 9
$ws.Range("F86").Value = "R4"

# Row 87
$ws.Range("D87").Value = "B"
$ws.Range("E87").Value = 11
$ws.Range("F87").Value = "R5"

# Row 88
$ws.Range("E88").Value = 44
$ws.Range("F88").Value = "R6"

# Row 89
$ws.Range("B89").Value = "F"
$ws.Range("C89").Value = "F"
$ws.Range("D89").Value = "D"
$ws.Range("E89").Value = 44
$ws.Range("F89").Value = "R6"

# Row 90
$ws.Range("B90").Value = "F"
$ws.Range("C90").Value = "F"
$ws.Range("D90").Value = "B"
$ws.Range("E90").Value = 11
$ws.Range("F90").Value = "R5"

# Row 91
$ws.Range("B91").Value = "F"
$ws.Range("C91").Value = "B"
$ws.Range("D91").Value = "B"
$ws.Range("E91").Value = 9
$ws.Range("F91").Value = "R4"

# Update view: scroll to new bottom area and change selection
$win = $excel.ActiveWindow
$win.ScrollRow = 73
$win.TopLeftCell = $ws.Range("A73")
$ws.Range("C94").Select() | Out-Null
